$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")
$ws.Rows("125:125").Delete()

$ws.Range("A150").Value = "PlacementStability"
$ws.Range("B150").Value = "cares\Placement.xlsx"
$ws.Range("C150").Value = "PlacementStability"
$ws.Range("D150").Value = 1

$ws.AutoFilterMode = $false
$ws.Range("A1:E138").AutoFilter()
